# Update the "Routing Master" sheet with the latest Engineering Item /
# Routing record (Item Number + Id) pulled in as part of the System Setup /
# CPQ test-plan refresh.
#
# Item Number (col B) -> "Pro-PEItem-28YRN"
# Id          (col D) -> "a345f000000uWBsAAM"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

$ws.Range("B2").Value = "Pro-PEItem-28YRN"
$ws.Range("D2").Value = "a345f000000uWBsAAM"
